# Remove the two duplicate phrase rows from the "list" sheet.

$wb = $excel.ActiveWorkbook
$list = $wb.Worksheets.Item("list")

# Delete the two duplicate-phrase rows (shifts everything below up by one
# each time). Delete bottom-up so the row numbers stay valid.
# Row 42: "[on a / It's been a] journey" (duplicate phrase)
# Row 41: "I've got to jump on another call" (duplicate of "I have to jump to another call")
$list.Rows.Item(42).EntireRow.Delete() | Out-Null
$list.Rows.Item(41).EntireRow.Delete() | Out-Null

# Leave the sheet's selection where it ended up after the rows shifted.
$list.Range("A45").Select() | Out-Null
